# Modified design of experiments function for training gaussian process model
#
# 1) Rename the "material-age-mrt" sheet to "material-age-velocity" (the
#    calibration experiment now sweeps on velocity instead of mean
#    residence time).
# 2) Add two new calibration runs (rows) to the "material-age-diameter"
#    results sheet.
# 3) Leave the renamed sheet active/selected, matching the author's saved
#    view (cell B1 selected, that sheet's tab active).

$wb = $excel.ActiveWorkbook

# --- 1. Rename sheet -------------------------------------------------
$wsVelocity = $wb.Worksheets.Item("material-age-mrt")
$wsVelocity.Name = "material-age-velocity"

# --- 2. Append new calibration rows to material-age-diameter ---------
$wsDiameter = $wb.Worksheets.Item("material-age-diameter")

$wsDiameter.Range("A3").Value = 20
$wsDiameter.Range("B3").Value = "wwmd"
$wsDiameter.Range("C3").Value = "kiosk only"
$wsDiameter.Range("D3").Value = -0.5126061925831822
$wsDiameter.Range("E3").Value = 0.001122609200955326
$wsDiameter.Range("F3").Value = 0.001275515361691412
$wsDiameter.Range("G3").Value = 15433.73804163933
$wsDiameter.Range("H3").Value = -0.1937379419510491
$wsDiameter.Range("I3").Value = -0.1072255400286369
$wsDiameter.Range("J3").Value = -0.03883206613208263
$wsDiameter.Range("K3").Value = -0.00001929727004422416
$wsDiameter.Range("L3").Value = -0.04997695755738575

$wsDiameter.Range("A4").Value = 20
$wsDiameter.Range("B4").Value = "wwmd"
$wsDiameter.Range("C4").Value = "kiosk only"
$wsDiameter.Range("D4").Value = -0.5126061925831822
$wsDiameter.Range("E4").Value = 0.001122609200955326
$wsDiameter.Range("F4").Value = 0.001275515361691412
$wsDiameter.Range("G4").Value = 15433.73804163933
$wsDiameter.Range("H4").Value = -0.1937379419510491
$wsDiameter.Range("I4").Value = -0.1072255400286369
$wsDiameter.Range("J4").Value = -0.03883206613208263
$wsDiameter.Range("K4").Value = -0.00001929727004422416
$wsDiameter.Range("L4").Value = -0.04997695755738575

# --- 3. Activate the renamed sheet and select B1 ----------------------
$wsVelocity.Activate()
$wsVelocity.Range("B1").Select()
